# Adds 5 new "Mac-Address"/device rows (reg center 10002, device ids
# 3000176-3000180) to the master-reg_center_device_h sheet, mirroring
# the existing rows' lang_code/is_active/cr_by/cr_dtimes/eff_dtimes
# values ("eng" / TRUE / "superadmin" / "now()" / "now()").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow    = 157
$regCenterId = 10002
$startDeviceId = 3000176
$rowCount    = 5

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $regCenterId
    $ws.Cells.Item($r, 2).Value = $startDeviceId + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the author's final view state: scrolled down with B157 selected.
$ws.Range("B157").Select()
try {
    $excel.ActiveWindow.ScrollRow = 152
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/limited hosts may not expose window scroll state; ignore.
}
